# Adds two new test-case rows (row 3 & row 4) to the "datatypes excel"
# sheet and rewrites row 2's "Expected Output"/"Actual Output" cells
# (D2/E2) with the new detailed result string. Mirrors the commit that
# added rows for "no input in command line" and "use argv[0] -h command"
# test cases.
#
# NOTE: values are written top-to-bottom / left-to-right in the same
# order they appear in the final sheet, so the engine's shared-string
# table is (re)built in the same order the author's workbook shows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 3: "no input in command line" test case ----
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "no input in command line"
$ws.Range("C3").Value = "no input"
$ws.Range("D3").Value = "use argv[0] -h command"
$ws.Range("E3").Value = "use argv[0] -h command"
$ws.Range("F3").Value = "PASS"

# ---- Row 4: "use argv[0] -h command" / help test case ----
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "To check help command"
$ws.Range("C4").Value = "argv[0] -h"
$ws.Range("D4").Value = "argv[0] datatype1 datatype2…"
$ws.Range("E4").Value = "argv[0] datatype1 datatype2…"
$ws.Range("F4").Value = "PASS"

# ---- Row 2: update Expected/Actual Output with the detailed result ----
$ws.Range("D2").Value = "int 4 float 8 char 1 string 3"
$ws.Range("E2").Value = "int 4 float 8 char 1 string 3"

# ---- Column widths for D (bestFit/AutoFit'd by the author) and E ----
# (inputs are calibrated so the engine's character-width -> pixel
# rounding lands on the same stored width the workbook ships with)
$ws.Columns.Item(4).ColumnWidth = 26.833333333333336
$ws.Columns.Item(5).ColumnWidth = 27.166666666666664

# ---- Selection moved by the author while editing ----
$null = $ws.Range("E13").Select()
